$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.130.88"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "1.795.93"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "314.16"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "0.5208"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("D8").Value = "0.3815"
$ws.Range("E8").Value = "  -3.71%  "
$ws.Range("D9").Value = "0.07904"
$ws.Range("E9").Value = "  -4.44%  "
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "1.099"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("D12").Value = "6.279"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "20.54"
$ws.Range("E14").Value = "  -3.17%  "
$ws.Range("D15").Value = "1.793.32"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("D16").Value = "7.227"
$ws.Range("E16").Value = "  -4.35%  "
$ws.Range("D17").Value = "93.15"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("E18").Value = "  -4.03%  "
$ws.Range("D19").Value = "0.06559"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "17.26"
$ws.Range("E21").Value = "  -3.28%  "
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("D23").Value = "28.184.64"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("E24").Value = "  -2.81%  "
$ws.Range("D25").Value = "2.269"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "160.86"
$ws.Range("E26").Value = "  +2.55%  "
$ws.Range("D27").Value = "20.43"
$ws.Range("E27").Value = "  -4.23%  "
$ws.Range("D28").Value = "2.000.26"
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("D29").Value = "2.334"
$ws.Range("E29").Value = "  -3.58%  "
$ws.Range("D30").Value = "122.94"
$ws.Range("E30").Value = "  -3.09%  "
$ws.Range("D31").Value = "0.1064"
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("E32").Value = "  -5.91%  "
$ws.Range("D33").Value = "3.675"
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").Value = "5.559"
$ws.Range("E34").Value = "  -3.96%  "
$ws.Range("D35").Value = "0.07307"
$ws.Range("E35").Value = "  +3.49%  "
$ws.Range("D36").Value = "12.22"
$ws.Range("E36").Value = "  +8.23%  "
$ws.Range("D37").Value = "0.02318"
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("D38").Value = "0.2139"
$ws.Range("E38").Value = "  -4.31%  "
$ws.Range("D39").Value = "8.688"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("D40").Value = "5.063"
$ws.Range("E40").Value = "  -3.93%  "
$ws.Range("D41").Value = "0.6144"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("D42").Value = "1.162"
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("D43").Value = "1.372"
$ws.Range("E43").Value = "  -2.25%  "
$ws.Range("D44").Value = "13.27"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("D45").Value = "3.781"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").Value = "0.5960"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "127.82"
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").Value = "1.231"
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("D49").Value = "1.917"
$ws.Range("E49").Value = "  -3.94%  "
$ws.Range("D50").Value = "0.06767"
$ws.Range("E50").Value = "  -2.42%  "
$ws.Range("D51").Value = "73.18"
$ws.Range("E51").Value = "  -1.57%  "
